$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: values in columns V ("c_Neu") and BB ("c_dash") for data rows
# 2-8 were placed in the wrong column; zero them out here.
$ws.Range("V2:V8").Value = 0
$ws.Range("BB2:BB8").Value = 0
